# default_sum_koef.xlsx - "updated default coef for archaeology"
#
# The element (column A) / coefficient (column B) pairs are re-ordered,
# two elements (Ta181, Nb93) are dropped and replaced by two new ones
# (P31, Co59), several coefficient values are refreshed, and column B is
# re-entered as text instead of numbers. Column A also picks up a
# left/vertically-centred alignment and a wider column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New element / coefficient pairs, in their final row order.
$rows = @(
    @{ A = "Na23";  B = "74.186" },
    @{ A = "Mg24";  B = "60.304" },
    @{ A = "Al27";  B = "52.925" },
    @{ A = "Si28";  B = "46.75" },
    @{ A = "P31";   B = "43.64" },
    @{ A = "K39";   B = "83.01" },
    @{ A = "Ca44";  B = "71.469" },
    @{ A = "Ti47";  B = "59.934" },
    @{ A = "Mn55";  B = "77.44" },
    @{ A = "Fe56";  B = "69.943" },
    @{ A = "Co59";  B = "78.65" },
    @{ A = "Cu63";  B = "79.88" },
    @{ A = "Sb121"; B = "83.53" },
    @{ A = "Pb208"; B = "92.83" },
    @{ A = "Sn118"; B = "78.77" }
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 1
    $ws.Range("A$r").Value = $rows[$i].A
    # Leading apostrophe forces text storage even though the value looks
    # numeric (matches the workbook's new t="s" coefficient cells).
    $ws.Range("B$r").Value = "'" + $rows[$i].B
    # Drop the quote-prefix flag the apostrophe entry leaves behind so the
    # cell keeps plain default formatting.
    $ws.Range("B$r").Style = "Normal"
}

# Column A: left-align + vertically centre every populated cell, and widen
# the column to fit the longest label.
$ws.Range("A1:A15").Style = "Normal"
$ws.Range("A1:A15").HorizontalAlignment = -4131
$ws.Range("A1:A15").VerticalAlignment = -4108
$ws.Columns("A").ColumnWidth = 9.65

# Selection moved to F13 in the saved file.
[void]$ws.Range("F13").Select()
